$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "64.665.39"
$ws.Cells.Item(2, 5).Value = "  +1.46%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.099.96"
$ws.Cells.Item(3, 5).Value = "  +1.51%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.03%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "564.69"
$ws.Cells.Item(5, 5).Value = "  +1.65%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "147.03"
$ws.Cells.Item(6, 5).Value = "  +2.12%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Cells.Item(7, 5).Value = "  +0.11%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "3.099.00"
$ws.Cells.Item(8, 5).Value = "  +1.54%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.521"
$ws.Cells.Item(9, 5).Value = "  +4.29%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +5.16%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -2.81%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +4.92%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000238"
$ws.Cells.Item(13, 5).Value = "  +5.58%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "36.50"
$ws.Cells.Item(14, 5).Value = "  +4.57%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.608.54"
$ws.Cells.Item(15, 5).Value = "  +1.05%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.653.93"
$ws.Cells.Item(16, 5).Value = "  +1.38%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.098.57"
$ws.Cells.Item(17, 5).Value = "  +1.47%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "7.00"
$ws.Cells.Item(18, 5).Value = "  +3.88%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.32%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "498.80"
$ws.Cells.Item(20, 5).Value = "  +5.49%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "14.69"
$ws.Cells.Item(21, 5).Value = "  +5.98%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "15.51"
$ws.Cells.Item(22, 5).Value = "  +15.44%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +5.60%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +2.73%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "83.82"
$ws.Cells.Item(25, 5).Value = "  +3.20%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.28%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +3.47%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "8.53"
$ws.Cells.Item(28, 5).Value = "  +6.44%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +3.45%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "27.72"
$ws.Cells.Item(30, 5).Value = "  +6.68%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.00"
$ws.Cells.Item(31, 5).Value = "  +0.01%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.66"
$ws.Cells.Item(32, 5).Value = "  +8.55%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +2.39%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.96"
$ws.Cells.Item(34, 5).Value = "  +7.72%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "6.51"
$ws.Cells.Item(35, 5).Value = "  +6.49%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "55.16"
$ws.Cells.Item(36, 5).Value = "  +1.08%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "456.42"
$ws.Cells.Item(37, 5).Value = "  -0.61%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0416"
$ws.Cells.Item(38, 5).Value = "  +3.45%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0837"
$ws.Cells.Item(39, 5).Value = "  +1.32%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.085.92"
$ws.Cells.Item(40, 5).Value = "  +4.77%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.84"
$ws.Cells.Item(41, 5).Value = "  -4.98%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.120"
$ws.Cells.Item(42, 5).Value = "  +5.41%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "8.47"
$ws.Cells.Item(43, 5).Value = "  +2.84%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.289"
$ws.Cells.Item(44, 5).Value = "  +12.24%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Fetch.AI"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.35"
$ws.Cells.Item(45, 5).Value = "  +10.60%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "28.58"
$ws.Cells.Item(46, 5).Value = "  +2.60%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -0.02%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +3.05%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0₃0540"
$ws.Cells.Item(49, 5).Value = "  +5.40%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +7.23%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "117.62"
$ws.Cells.Item(51, 5).Value = "  -1.48%  "
